$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -414685708118954304.0
$ws.Range("C2").Value = -172000681283411.0
$ws.Range("D2").Value = 101689501379444527333376.0
$ws.Range("E2").Value = 9834185715030458368.0
$ws.Range("B3").Value = -414103102252363712.0
$ws.Range("C3").Value = -171759823244871.0
$ws.Range("D3").Value = 101552918427403292770304.0
$ws.Range("E3").Value = 9820316061320046592.0
$ws.Range("B4").Value = -413520536541290240.0
$ws.Range("C4").Value = -171519150536391.0
$ws.Range("D4").Value = 101416381568627412828160.0
$ws.Range("E4").Value = 9806446419267573760.0
$ws.Range("B5").Value = -412938002474007424.0
$ws.Range("C5").Value = -171278081104149.0
$ws.Range("D5").Value = 101279851876151065575424.0
$ws.Range("E5").Value = 9792576787378438144.0
$ws.Range("B6").Value = -412355499551199040.0
$ws.Range("C6").Value = -171036892359019.0
$ws.Range("D6").Value = 101143351766159684796416.0
$ws.Range("E6").Value = 9778707166723219456.0
$ws.Range("B7").Value = -411773001050911616.0
$ws.Range("C7").Value = -170795455428471.0
$ws.Range("D7").Value = 101006886204998412140544.0
$ws.Range("E7").Value = 9764837546392616960.0
$ws.Range("B8").Value = -411190488995108416.0
$ws.Range("C8").Value = -170553582028346.0
$ws.Range("D8").Value = 100870207038307352707072.0
$ws.Range("E8").Value = 9750967922969516032.0
$ws.Range("B9").Value = -410607972748840896.0
$ws.Range("C9").Value = -170312488773587.0
$ws.Range("D9").Value = 100733437300810723622912.0
$ws.Range("E9").Value = 9737098307609286656.0
$ws.Range("B10").Value = -410025502771213248.0
$ws.Range("C10").Value = -170071166797970.0
$ws.Range("D10").Value = 100596939768077198819328.0
$ws.Range("E10").Value = 9723228697136887808.0
$ws.Range("B11").Value = -409443111933909696.0
$ws.Range("C11").Value = -169828242700886.0
$ws.Range("D11").Value = 100460307232916055785472.0
$ws.Range("E11").Value = 9709359113879023616.0
$ws.Range("B12").Value = -408860665211106304.0
$ws.Range("C12").Value = -169586493747859.0
$ws.Range("D12").Value = 100323684236277323923456.0
$ws.Range("E12").Value = 9695489511912505344.0
$ws.Range("B13").Value = -408278231447801536.0
$ws.Range("C13").Value = -169344880724335.0
$ws.Range("D13").Value = 100187138660363090788352.0
$ws.Range("E13").Value = 9681619910661543936.0
$ws.Range("B14").Value = -407695781375954048.0
$ws.Range("C14").Value = -169103061725936.0
$ws.Range("D14").Value = 100050871864227541286912.0
$ws.Range("E14").Value = 9667750307175444480.0
$ws.Range("B15").Value = -407113316052572736.0
$ws.Range("C15").Value = -168859078342954.0
$ws.Range("D15").Value = 99914434130193806262272.0
$ws.Range("E15").Value = 9653880698307862528.0
$ws.Range("B16").Value = -406530915123067712.0
$ws.Range("C16").Value = -168615040441461.0
$ws.Range("D16").Value = 99778069783030146269184.0
$ws.Range("E16").Value = 9640011110852474880.0
$ws.Range("B17").Value = -405948474345990720.0
$ws.Range("C17").Value = -168369030857937.0
$ws.Range("D17").Value = 99641636674204011069440.0
$ws.Range("E17").Value = 9626141511763718144.0
$ws.Range("B18").Value = -405365953959028544.0
$ws.Range("C18").Value = -168123491646151.0
$ws.Range("D18").Value = 99505291342404393107456.0
$ws.Range("E18").Value = 9612271884470534144.0
$ws.Range("B19").Value = -404783557479569600.0
$ws.Range("C19").Value = -167879544878563.0
$ws.Range("D19").Value = 99369010004353037107200.0
$ws.Range("E19").Value = 9598402297237733376.0
$ws.Range("B20").Value = -404201098444760384.0
$ws.Range("C20").Value = -167639986031700.0
$ws.Range("D20").Value = 99232636886439581188096.0
$ws.Range("E20").Value = 9584532688124192768.0
$ws.Range("B21").Value = -403618762651959936.0
$ws.Range("C21").Value = -167395151366887.0
$ws.Range("D21").Value = 99096348829084645588992.0
$ws.Range("E21").Value = 9570663122144473088.0
$ws.Range("B22").Value = -403036277523728448.0
$ws.Range("C22").Value = -167149343184602.0
$ws.Range("D22").Value = 98959875858994510692352.0
$ws.Range("E22").Value = 9556793512277878784.0
$ws.Range("B23").Value = -402453963208102720.0
$ws.Range("C23").Value = -166905320670277.0
$ws.Range("D23").Value = 98823958405789708713984.0
$ws.Range("E23").Value = 9542923951324366848.0
$ws.Range("B24").Value = -401871608876554112.0
$ws.Range("C24").Value = -166661154095644.0
$ws.Range("D24").Value = 98687937481902978498560.0
$ws.Range("E24").Value = 9529054378289557504.0
$ws.Range("B25").Value = -401289233000004160.0
$ws.Range("C25").Value = -166417736963711.0
$ws.Range("D25").Value = 98551684711607682727936.0
$ws.Range("E25").Value = 9515184800921980928.0
$ws.Range("B26").Value = -400706867269414720.0
$ws.Range("C26").Value = -166176215139846.0
$ws.Range("D26").Value = 98415446032347945762816.0
$ws.Range("E26").Value = 9501315223276883968.0
$ws.Range("B27").Value = -400124457257692672.0
$ws.Range("C27").Value = -165932495624609.0
$ws.Range("D27").Value = 98279397793893934170112.0
$ws.Range("E27").Value = 9487445630468073472.0
$ws.Range("B28").Value = -399542057905635392.0
$ws.Range("C28").Value = -165686731365111.0
$ws.Range("D28").Value = 98143285890155106271232.0
$ws.Range("E28").Value = 9473576036040912896.0
